$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("I8").Value = "type: blog`nwidth: 2`nheight: 1`nser: 128"
